$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the "Basic ui test" row (row 5) as done in the Status column
$ws.Range("B5").Value = "done"

# Add a new to-do item in row 11
$ws.Range("A11").Value = "Implement CI/CD"

# Move the active selection to the newly added cell
$ws.Range("A11").Select()
